$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed from the source data ("RM 232" and "SC 92").
# These are located by matching column A labels rather than hard-coded row numbers,
# since row positions shift as rows are removed.
$rm232 = $ws.Range("A1:A35").Find("RM 232")
if ($rm232 -ne $null) {
    $ws.Rows.Item($rm232.Row).Delete()
}

$sc92 = $ws.Range("A1:A35").Find("SC 92")
if ($sc92 -ne $null) {
    $ws.Rows.Item($sc92.Row).Delete()
}

# Apply the individual cell-value edits (addressed by the final, post-deletion row numbers).
$ws.Range("F5").ClearContents()
$ws.Range("F11").Value = 17.65
$ws.Range("E19").Value = -6.5
$ws.Range("F19").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("F25").Value = 16.6
$ws.Range("E27").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("E33").Value = -10.7
